# Fix mis-assigned "garant" rows: for several katedra+zkratka groups the
# A (katedra), E (jednotekCviceni) and F (jednotkaCviceni) values had been
# written to the wrong one of the two sibling rows. Swap them back.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs (1-based worksheet rows) whose A/E/F columns must be swapped.
$pairs = @(
    @(29, 30),
    @(51, 52),
    @(66, 67),
    @(69, 70),
    @(88, 89)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $a1 = $ws.Cells.Item($r1, 1).Value2
    $e1 = $ws.Cells.Item($r1, 5).Value2
    $f1 = $ws.Cells.Item($r1, 6).Value2

    $a2 = $ws.Cells.Item($r2, 1).Value2
    $e2 = $ws.Cells.Item($r2, 5).Value2
    $f2 = $ws.Cells.Item($r2, 6).Value2

    $ws.Cells.Item($r1, 1).Value = $a2
    $ws.Cells.Item($r1, 5).Value = $e2
    $ws.Cells.Item($r1, 6).Value = $f2

    $ws.Cells.Item($r2, 1).Value = $a1
    $ws.Cells.Item($r2, 5).Value = $e1
    $ws.Cells.Item($r2, 6).Value = $f1
}
